$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 110: explain the status bits for the REMOTE UPDATE status register
# ---------------------------------------------------------------------------
$ws.Range("D110").Value = "(0):RU busy ; (1): EPCQ busy ; (2) EPCQ done ; (3) FIFO rd  empty ; (4) FIFO wr full"

# ---------------------------------------------------------------------------
# 2. Rows 113 / 114: label the two new read-only EPCQ data registers
# ---------------------------------------------------------------------------
$ws.Range("C113").Value = "REMOTE UPDATE RD ONLY-> epcq data low"
$ws.Range("C114").Value = "REMOTE UPDATE RD ONLY-> epcq data high"

# ---------------------------------------------------------------------------
# 3. Relocate the "clock select" register (currently row 127) down to row
#    131 so that four new rows (113, 114, 127, 128) can describe the new
#    EPCQ-read-RAM registers that were inserted ahead of it.
# ---------------------------------------------------------------------------
$ws.Range("C127:F127").Cut($ws.Range("C131:F131"))

# remove the left-over (now blank) formatted cells in row 127
$ws.Range("D127:F127").Clear()

# ---------------------------------------------------------------------------
# 4. New register rows describing the EPCQ read-RAM address/enable/clk
# ---------------------------------------------------------------------------
$ws.Range("C127").Value = "REMOTE UPDATE(10) -> EPCQ read RAM addr (lower 12 bits)"

$ws.Range("C128").Value = "REMOTE UPDATE(11) -> EPCQ read RAM clk/enable"
$ws.Range("D128").Value = "bit 0 : enable, bit 1 : clk"

# ---------------------------------------------------------------------------
# 5. Update the view so the sheet selects where the author left off
# ---------------------------------------------------------------------------
$ws.Range("C128").Select()
